$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force affected Price (D) cells to Text format so numeric-looking strings
# (e.g. "605.64") are preserved as text instead of being parsed as numbers,
# matching the original inline-string cell type.
$dCells = @('D2','D3','D5','D6','D8','D9','D13','D14','D16','D17','D18','D19','D20','D21','D23','D24','D25','D26','D29','D36','D37','D38','D39','D41','D43','D44','D45','D46','D47','D48','D49','D50')
foreach ($addr in $dCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '67.872.25'
$ws.Range('E2').Value = '  +1.90%  '

$ws.Range('D3').Value = '2.625.59'
$ws.Range('E3').Value = '  +2.22%  '

$ws.Range('E4').Value = '  -0.16%  '

$ws.Range('D5').Value = '605.64'
$ws.Range('E5').Value = '  +2.78%  '

$ws.Range('D6').Value = '154.81'
$ws.Range('E6').Value = '  +1.19%  '

$ws.Range('E7').Value = '  -0.02%  '

$ws.Range('D8').Value = '0.551'
$ws.Range('E8').Value = '  +2.46%  '

$ws.Range('D9').Value = '2.622.10'
$ws.Range('E9').Value = '  +2.01%  '

$ws.Range('E10').Value = '  +14.32%  '

$ws.Range('E11').Value = '  +0.88%  '

$ws.Range('E12').Value = '  +1.82%  '

$ws.Range('D13').Value = '0.356'
$ws.Range('E13').Value = '  +0.95%  '

$ws.Range('D14').Value = '27.90'
$ws.Range('E14').Value = '  +0.32%  '

$ws.Range('E15').Value = '  +5.87%  '

$ws.Range('D16').Value = '3.096.41'
$ws.Range('E16').Value = '  +1.53%  '

$ws.Range('D17').Value = '67.862.20'
$ws.Range('E17').Value = '  +1.80%  '

$ws.Range('D18').Value = '2.618.50'
$ws.Range('E18').Value = '  +1.83%  '

$ws.Range('D19').Value = '11.26'
$ws.Range('E19').Value = '  +0.93%  '

$ws.Range('D20').Value = '366.59'
$ws.Range('E20').Value = '  +4.25%  '

$ws.Range('D21').Value = '7.69'
$ws.Range('E21').Value = '  -0.08%  '

$ws.Range('E22').Value = '  -0.03%  '

$ws.Range('D23').Value = '2.10'
$ws.Range('E23').Value = '  +3.67%  '

$ws.Range('D24').Value = '1.00'
$ws.Range('E24').Value = '  -0.04%  '

$ws.Range('D25').Value = '70.46'
$ws.Range('E25').Value = '  +4.97%  '

$ws.Range('D26').Value = '9.98'
$ws.Range('E26').Value = '  -2.83%  '

$ws.Range('E27').Value = '  +3.97%  '

$ws.Range('D29').Value = '582.30'
$ws.Range('E29').Value = '  -1.22%  '

$ws.Range('E30').Value = '  -0.05%  '

$ws.Range('E31').Value = '  +0.31%  '

$ws.Range('E32').Value = '  +0.07%  '

$ws.Range('E33').Value = '  +2.25%  '

$ws.Range('E34').Value = '  -0.63%  '

$ws.Range('E35').Value = '  +0.00%  '

$ws.Range('D36').Value = '1.54'
$ws.Range('E36').Value = '  -1.29%  '

$ws.Range('D37').Value = '4.98'
$ws.Range('E37').Value = '  +0.71%  '

$ws.Range('D38').Value = '157.70'
$ws.Range('E38').Value = '  +2.75%  '

$ws.Range('D39').Value = '19.49'
$ws.Range('E39').Value = '  +2.55%  '

$ws.Range('E40').Value = '  +1.49%  '

$ws.Range('D41').Value = '5.40'
$ws.Range('E41').Value = '  +0.11%  '

$ws.Range('E42').Value = '  +4.89%  '

$ws.Range('D43').Value = '2.65'
$ws.Range('E43').Value = '  +2.30%  '

$ws.Range('D44').Value = '41.15'
$ws.Range('E44').Value = '  -0.64%  '

$ws.Range('B45').Value = 'WhiteBITCoin'
$ws.Range('C45').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D45').Value = '16.44'
$ws.Range('E45').Value = '  +0.22%  '

$ws.Range('B46').Value = 'USDe'
$ws.Range('C46').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D46').Value = '1.00'
$ws.Range('E46').Value = '  +0.05%  '

$ws.Range('D47').Value = '157.53'
$ws.Range('E47').Value = '  +2.37%  '

$ws.Range('D48').Value = '0.0₆0291'
$ws.Range('E48').Value = '  -4.48%  '

$ws.Range('D49').Value = '3.77'
$ws.Range('E49').Value = '  +1.28%  '

$ws.Range('D50').Value = '21.05'
$ws.Range('E50').Value = '  +0.75%  '

$ws.Range('E51').Value = '  +2.50%  '

# Restore default cell style on the Price cells (removes the temporary
# Text number format applied above) so only the values/content differ.
foreach ($addr in $dCells) {
    $ws.Range($addr).Style = "Normal"
}
